# MEv1 workbook update
# - Refresh "Planilla" expense and dependent Total formula result
# - Rename "Lista por Pagar:" section to "Lista de Egresos:" and add a
#   parallel "Por Hacer:" to-do list
# - Add a few new items to the existing "Accesorios" and "Vestimenta"
#   shopping sub-lists

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Planilla amount changed; the Total (C8 = SUM(C4:C7)) recalculates on its own
$ws.Range("C4").Value = 818.37

# New "Accesorios" entries (column R, under header R14 "Accesorios")
$ws.Range("R16").Value = "Desodorantes"
$ws.Range("R17").Value = "Perfumes"

# New "Vestimenta" entries (column P, under header P14 "Vestimenta")
$ws.Range("P20").Value = "Medias"
$ws.Range("P21").Value = "Short"

# Section header F24 renamed from "Lista por Pagar:" to "Lista de Egresos:"
$ws.Range("F24").Value = "Lista de Egresos:"

# New "Por Hacer:" list alongside the "Lista de Egresos:" column,
# bold like the other section headers (F12, F24)
$ws.Range("H24").Value = "Por Hacer:"
$ws.Range("H24").Font.Bold = $true

$ws.Range("H26").Value = "Arreglar pluma"
$ws.Range("H27").Value = "Switch de cocina"
$ws.Range("H28").Value = "Sacar Lentes"

# Move the active selection, matching the author's last cursor position
$ws.Range("P22").Select() | Out-Null

Write-Output "edit applied"
